$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 0.8732805848121643
$ws.Range("B1").Value = 2.994909286499023
$ws.Range("C1").Value = 4.522442817687988
$ws.Range("D1").Value = 2.908376216888428
$ws.Range("E1").Value = 1.418692469596863
